$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the total marks error: update Right/Wrong counts and the Max summary text.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -16
$ws.Range("E12").Value = "44 / 112"
